$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 507 (old rows 507..619 shift down to 509..621)
$ws.Rows("507:508").Insert()

# New row 507
$ws.Cells.Item(507, 1).Value = 6
$ws.Cells.Item(507, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(507, 3).Value = "Metropolitana"
$ws.Cells.Item(507, 4).Value = 44951
$ws.Cells.Item(507, 5).Value = 13
$ws.Cells.Item(507, 6).Value = 100112030
$ws.Cells.Item(507, 7).Value = "Poroto granado"
$ws.Cells.Item(507, 8).Value = "Sin especificar"
$ws.Cells.Item(507, 9).Value = "Primera"
$ws.Cells.Item(507, 10).Value = 510
$ws.Cells.Item(507, 11).Value = 35000
$ws.Cells.Item(507, 12).Value = 37000
$ws.Cells.Item(507, 13).Value = 36098
$ws.Cells.Item(507, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(507, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(507, 16).Value = 1444
$ws.Cells.Item(507, 17).Value = 25
$ws.Cells.Item(507, 18).Value = "Hortaliza"

# New row 508
$ws.Cells.Item(508, 1).Value = 6
$ws.Cells.Item(508, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(508, 3).Value = "Metropolitana"
$ws.Cells.Item(508, 4).Value = 44951
$ws.Cells.Item(508, 5).Value = 13
$ws.Cells.Item(508, 6).Value = 100112030
$ws.Cells.Item(508, 7).Value = "Poroto granado"
$ws.Cells.Item(508, 8).Value = "Sin especificar"
$ws.Cells.Item(508, 9).Value = "Primera"
$ws.Cells.Item(508, 10).Value = 600
$ws.Cells.Item(508, 11).Value = 35000
$ws.Cells.Item(508, 12).Value = 38000
$ws.Cells.Item(508, 13).Value = 36850
$ws.Cells.Item(508, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(508, 15).Value = "Región Metropolitana"
$ws.Cells.Item(508, 16).Value = 1474
$ws.Cells.Item(508, 17).Value = 25
$ws.Cells.Item(508, 18).Value = "Hortaliza"

Write-Output "done"
